# Apply the new Mega-Sena draw rows (contests 2895-2914) to the "MEGA SENA" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEGA SENA")

$data = @(
    @(2895, 9, 11, 44, 51, 52, 56),
    @(2896, 8, 9, 12, 16, 43, 53),
    @(2897, 1, 6, 24, 27, 28, 57),
    @(2898, 14, 23, 30, 32, 38, 48),
    @(2899, 10, 22, 28, 42, 44, 51),
    @(2900, 33, 50, 54, 55, 59, 60),
    @(2901, 2, 20, 28, 38, 44, 47),
    @(2902, 8, 21, 22, 42, 45, 48),
    @(2903, 20, 24, 27, 46, 50, 54),
    @(2904, 2, 37, 38, 46, 52, 55),
    @(2905, 4, 17, 18, 26, 43, 52),
    @(2906, 17, 33, 37, 41, 46, 49),
    @(2907, 30, 33, 42, 44, 52, 56),
    @(2908, 20, 35, 36, 37, 38, 50),
    @(2909, 8, 21, 31, 41, 53, 58),
    @(2910, 3, 4, 11, 15, 28, 29),
    @(2911, 23, 27, 32, 54, 56, 59),
    @(2912, 9, 25, 37, 41, 51, 59),
    @(2913, 17, 21, 34, 52, 55, 60),
    @(2914, 18, 25, 35, 40, 46, 47)
)

$startRow = 352
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$endRow = $startRow + $data.Count - 1

# Move the selection to the newly added range, matching the saved view state.
$ws.Range("A$startRow" + ":G$endRow").Select() | Out-Null
